# Replace the old nickname/chat text "dikakokodikakoko" with "akudikakoko"
# on the Chat sheet. This removes the now-unused shared string
# "dikakokodikakoko" from the workbook's shared string table (causing every
# later shared-string index to shift down by one) and appends the brand new
# string "akudikakoko" at the end of the table.
$wb = $excel.ActiveWorkbook

$chat = $wb.Worksheets.Item("Chat")
$chat.Cells.Replace("dikakokodikakoko", "akudikakoko")

# Move the active sheet / selection from Login_Alert_Popup to Chat.
# First restore the left-behind selection on Login_Alert_Popup (it keeps
# its own remembered selection even once it's no longer the active sheet).
$loginAlertPopup = $wb.Worksheets.Item("Login_Alert_Popup")
$loginAlertPopup.Activate()
$loginAlertPopup.Range("I5").Select()

# Now activate Chat and move its selection, leaving it as the active sheet.
$chat.Activate()
$chat.Range("E9").Select()
